$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (year 2025) metrics per latest data refresh
$ws.Range("C8").Value = 1254
$ws.Range("D8").Value = 202
$ws.Range("E8").Value = 1052
$ws.Range("F8").Value = 8.285479901558656
$ws.Range("G8").Value = 83.89154704944178
$ws.Range("H8").Value = 16.10845295055821
